$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.669548153877258
$ws.Range("B1").Value = 2.064101696014404
$ws.Range("C1").Value = 2.224788665771484
$ws.Range("D1").Value = 2.56936502456665
$ws.Range("E1").Value = 3.316083669662476
